$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add experimentDesign column (D) for all data rows: "90minuteInduction"
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 4).Value = "90minuteInduction"
}

# Update harvester column (B) for all data rows: "Retrofitted_0648" -> "S.GISH"
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 2).Value = "S.GISH"
}

# Update the active selection to match the saved workbook state (D19)
$ws.Range("D19").Select()
